$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 148.82857
$ws.Range("I33").Value = 150.55882
$ws.Range("K33").Value = 150.55882
$ws.Range("M33").Value = 78.44118
$ws.Range("H43").Value = 3129.9
$ws.Range("I43").Value = 3074.5
$ws.Range("J43").Value = 3166.8333
$ws.Range("K43").Value = 3074.5
$ws.Range("L43").Value = 3166.8333
$ws.Range("M43").Value = -3005.5
$ws.Range("N43").Value = -3304.8333
$ws.Range("H54").Value = 1999
$ws.Range("I54").Value = 1999
$ws.Range("K54").Value = 1999
$ws.Range("M54").Value = -1513
$ws.Range("H62").Value = 10531.077
$ws.Range("I62").Value = 11655.8
$ws.Range("K62").Value = 11655.8
$ws.Range("M62").Value = -11031.8
$ws.Range("H65").Value = 10531.077
$ws.Range("I65").Value = 11655.8
$ws.Range("K65").Value = 58279
$ws.Range("M65").Value = -55159
$ws.Range("I80").Value = 1001
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3003
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2005
$ws.Range("N80").ClearContents()
$ws.Range("I83").Value = 1001
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9009
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4017
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 3595
$ws.Range("I86").Value = 3600
$ws.Range("K86").Value = 3600
$ws.Range("M86").Value = -2477
$ws.Range("H89").Value = 3595
$ws.Range("I89").Value = 3600
$ws.Range("K89").Value = 18000
$ws.Range("M89").Value = -12384
$ws.Range("H127").Value = 1104.375
$ws.Range("I127").Value = 1047.8572
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 3143.5716
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 1816.4284
$ws.Range("N127").Value = -14420
$ws.Range("H137").Value = 2353.923
$ws.Range("I137").Value = 1765.1666
$ws.Range("K137").Value = 5295.4998
$ws.Range("M137").Value = -2745.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3741.6667
$ws.Range("I2").Value = 2611.818
$ws.Range("K2").Value = 2611.818
$ws.Range("M2").Value = -2498.818
$ws.Range("H116").Value = 3741.6667
$ws.Range("I116").Value = 2611.818
$ws.Range("K116").Value = 2611.818
$ws.Range("M116").Value = -317.8180000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3741.6667
$ws.Range("I3").Value = 2611.818
$ws.Range("K3").Value = 2611.818
$ws.Range("M3").Value = -2497.818
$ws.Range("H134").Value = 2445.7837
$ws.Range("I134").Value = 2508.9092
$ws.Range("J134").Value = 1925
$ws.Range("K134").Value = 7526.7276
$ws.Range("L134").Value = 5775
$ws.Range("M134").Value = -4991.7276
$ws.Range("N134").Value = -10845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 18010
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H25").Value = 9997.5
$ws.Range("I25").Value = 9997.5
$ws.Range("K25").Value = 9997.5
$ws.Range("M25").Value = -9823.5
$ws.Range("H27").Value = 18010
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 200002820
$ws.Range("I112").Value = 333334700
$ws.Range("K112").Value = 1000004100
$ws.Range("M112").Value = -1000002992
$ws.Range("H128").Value = 540208.3
$ws.Range("I128").Value = 540208.3
$ws.Range("K128").Value = 1620624.9
$ws.Range("M128").Value = -1615644.9
$ws.Range("H131").Value = 71431700
$ws.Range("I131").Value = 125000616
$ws.Range("J131").Value = 6480
$ws.Range("K131").Value = 375001848
$ws.Range("L131").Value = 19440
$ws.Range("M131").Value = -374996808
$ws.Range("N131").Value = -29520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 25635
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 33452.5
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 33452.5
$ws.Range("M22").Value = -9471
$ws.Range("N22").Value = -34510.5
$ws.Range("H31").Value = 9270
$ws.Range("I31").Value = 8125
$ws.Range("J31").Value = 13850
$ws.Range("K31").Value = 8125
$ws.Range("L31").Value = 13850
$ws.Range("M31").Value = -7833
$ws.Range("N31").Value = -14434
$ws.Range("H37").Value = 9270
$ws.Range("I37").Value = 8125
$ws.Range("J37").Value = 13850
$ws.Range("K37").Value = 8125
$ws.Range("L37").Value = 13850
$ws.Range("M37").Value = -7848
$ws.Range("N37").Value = -14404
$ws.Range("H113").Value = 3297.389
$ws.Range("I113").Value = 2740.8125
$ws.Range("K113").Value = 2740.8125
$ws.Range("M113").Value = -570.8125
$ws.Range("H126").Value = 4245.364
$ws.Range("I126").Value = 2814.2856
$ws.Range("K126").Value = 8442.856800000001
$ws.Range("M126").Value = -5972.856800000001
$ws.Range("H132").Value = 3896.8447
$ws.Range("I132").Value = 3676.5483
$ws.Range("J132").Value = 4149.778
$ws.Range("K132").Value = 11029.6449
$ws.Range("L132").Value = 12449.334
$ws.Range("M132").Value = -8499.644899999999
$ws.Range("N132").Value = -17509.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 13752.883
$ws.Range("I46").Value = 2472
$ws.Range("K46").Value = 2472
$ws.Range("M46").Value = -2284
$ws.Range("H61").Value = 80884.92
$ws.Range("I61").Value = 103200.5
$ws.Range("K61").Value = 103200.5
$ws.Range("M61").Value = -102998.5
$ws.Range("H113").Value = 80884.92
$ws.Range("I113").Value = 103200.5
$ws.Range("K113").Value = 103200.5
$ws.Range("M113").Value = -101030.5
$ws.Range("H116").Value = 69000
$ws.Range("J116").Value = 69000
$ws.Range("L116").Value = 69000
$ws.Range("N116").Value = -78178
$ws.Range("H120").Value = 67999
$ws.Range("J120").Value = 67999
$ws.Range("L120").Value = 67999
$ws.Range("N120").Value = -77675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 17500
$ws.Range("J11").Value = 17500
$ws.Range("L11").Value = 17500
$ws.Range("N11").Value = -17784
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H106").Value = 44980.668
$ws.Range("J106").Value = 44980.668
$ws.Range("L106").Value = 44980.668
$ws.Range("N106").Value = -47504.668
$ws.Range("H113").Value = 526.3913
$ws.Range("I113").Value = 514.1
$ws.Range("J113").Value = 608.3333
$ws.Range("K113").Value = 1542.3
$ws.Range("L113").Value = 1824.9999
$ws.Range("M113").Value = 627.6999999999998
$ws.Range("N113").Value = -6164.9999
